# Auto-generated edit script applying numeric updates to Titan_Profits sheets
# per the commit diff (8 worksheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Range("H16").Value = 7000
$ws.Range("J16").Value = 7000
$ws.Range("L16").Value = 7000
$ws.Range("N16").Value = -7460
# Row 137
$ws.Range("H137").Value = 142861420
$ws.Range("I137").Value = 142861420
$ws.Range("K137").Value = 428584260
$ws.Range("M137").Value = -428581710
# Row 138
$ws.Range("H138").Value = 5983192.5
$ws.Range("I138").Value = 1197148.1
$ws.Range("K138").Value = 3591444.3
$ws.Range("M138").Value = -3586304.3

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
# Row 33
$ws.Range("H33").Value = 10600
$ws.Range("I33").Value = 7000
$ws.Range("K33").Value = 7000
$ws.Range("M33").Value = -6671
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
# Row 61
$ws.Range("H61").Value = 3130.2273
$ws.Range("I61").Value = 2292.5334
$ws.Range("K61").Value = 2292.5334
$ws.Range("M61").Value = -2080.5334
# Row 74
$ws.Range("H74").Value = 5889.5557
$ws.Range("I74").Value = 1653.8948
$ws.Range("J74").Value = 15949.25
$ws.Range("K74").Value = 1653.8948
$ws.Range("L74").Value = 15949.25
$ws.Range("M74").Value = -779.8948
$ws.Range("N74").Value = -17697.25
# Row 77
$ws.Range("H77").Value = 5889.5557
$ws.Range("I77").Value = 1653.8948
$ws.Range("J77").Value = 15949.25
$ws.Range("K77").Value = 8269.474
$ws.Range("L77").Value = 79746.25
$ws.Range("M77").Value = -3901.474
$ws.Range("N77").Value = -88482.25
# Row 136
$ws.Range("H136").Value = 3130.2273
$ws.Range("I136").Value = 2292.5334
$ws.Range("K136").Value = 6877.600199999999
$ws.Range("M136").Value = -4327.600199999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 14
$ws.Range("H14").Value = 50000
$ws.Range("I14").Value = 50000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 50000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -49828
$ws.Range("N14").ClearContents()
# Row 17
$ws.Range("H17").Value = 3483.3333
$ws.Range("I17").Value = 450
$ws.Range("J17").Value = 5000
$ws.Range("K17").Value = 450
$ws.Range("L17").Value = 5000
$ws.Range("M17").Value = -278
$ws.Range("N17").Value = -5344
# Row 18
$ws.Range("H18").Value = 50000
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
# Row 22
$ws.Range("H22").Value = 174.93333
$ws.Range("I22").Value = 158
$ws.Range("K22").Value = 158
$ws.Range("M22").Value = 15
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
# Row 134
$ws.Range("H134").Value = 2821.7666
$ws.Range("I134").Value = 1776.7142
$ws.Range("J134").Value = 5260.222
$ws.Range("K134").Value = 5330.142599999999
$ws.Range("L134").Value = 15780.666
$ws.Range("M134").Value = -2795.142599999999
$ws.Range("N134").Value = -20850.666

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 403.5
$ws.Range("I10").Value = 403.5
$ws.Range("K10").Value = 403.5
$ws.Range("M10").Value = -264.5
# Row 13
$ws.Range("H13").Value = 205080
$ws.Range("J13").Value = 341666.66
$ws.Range("L13").Value = 341666.66
$ws.Range("N13").Value = -341944.66
# Row 15
$ws.Range("H15").Value = 1000
$ws.Range("I15").Value = 1000
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -830
$ws.Range("N15").ClearContents()
# Row 22
$ws.Range("H22").Value = 139.5
$ws.Range("J22").Value = 180
$ws.Range("L22").Value = 180
$ws.Range("N22").Value = -880
# Row 31
$ws.Range("H31").Value = 1155.238
$ws.Range("I31").Value = 1155.238
$ws.Range("K31").Value = 1155.238
$ws.Range("M31").Value = -860.2380000000001
# Row 32
$ws.Range("H32").Value = 17801.4
$ws.Range("I32").Value = 19751.75
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 19751.75
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -19435.75
$ws.Range("N32").Value = -10632
# Row 34
$ws.Range("H34").Value = 1155.238
$ws.Range("I34").Value = 1155.238
$ws.Range("K34").Value = 1155.238
$ws.Range("M34").Value = -953.2380000000001
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
# Row 58
$ws.Range("H58").Value = 1633.1428
$ws.Range("I58").Value = 1192.36
$ws.Range("J58").Value = 2735.1
$ws.Range("K58").Value = 1192.36
$ws.Range("L58").Value = 2735.1
$ws.Range("M58").Value = -989.3599999999999
$ws.Range("N58").Value = -3141.1
# Row 132
$ws.Range("H132").Value = 2953.9092
$ws.Range("I132").Value = 1876.2858
$ws.Range("J132").Value = 4839.75
$ws.Range("K132").Value = 5628.857400000001
$ws.Range("L132").Value = 14519.25
$ws.Range("M132").Value = -3098.857400000001
$ws.Range("N132").Value = -19579.25
# Row 134
$ws.Range("H134").Value = 2293.0557
$ws.Range("I134").Value = 1188.8889
$ws.Range("J134").Value = 5605.5557
$ws.Range("K134").Value = 3566.6667
$ws.Range("L134").Value = 16816.6671
$ws.Range("M134").Value = -1031.6667
$ws.Range("N134").Value = -21886.6671
# Row 136
$ws.Range("H136").Value = 1633.1428
$ws.Range("I136").Value = 1192.36
$ws.Range("J136").Value = 2735.1
$ws.Range("K136").Value = 3577.08
$ws.Range("L136").Value = 8205.299999999999
$ws.Range("M136").Value = -1027.08
$ws.Range("N136").Value = -13305.3

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1372.8494
$ws.Range("J131").Value = 1527.5646
$ws.Range("L131").Value = 4582.6938
$ws.Range("N131").Value = -14662.6938

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
# Row 17
$ws.Range("H17").Value = 50000
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
# Row 21
$ws.Range("H21").Value = 3959.8
$ws.Range("I21").Value = 98
$ws.Range("J21").Value = 4388.8887
$ws.Range("K21").Value = 98
$ws.Range("L21").Value = 4388.8887
$ws.Range("M21").Value = 75
$ws.Range("N21").Value = -4734.8887
# Row 30
$ws.Range("H30").Value = 3959.8
$ws.Range("I30").Value = 98
$ws.Range("J30").Value = 4388.8887
$ws.Range("K30").Value = 98
$ws.Range("L30").Value = 4388.8887
$ws.Range("M30").Value = 7
$ws.Range("N30").Value = -4598.8887
# Row 126
$ws.Range("H126").Value = 2708.24
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2885.3
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 8655.900000000001
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -13595.9
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 135
$ws.Range("H135").Value = 250015000
$ws.Range("J135").Value = 250015000
$ws.Range("L135").Value = 250015000
$ws.Range("N135").Value = -250025140
# Row 137
$ws.Range("H137").Value = 64999.5
$ws.Range("J137").Value = 64999.5
$ws.Range("L137").Value = 64999.5
$ws.Range("N137").Value = -75199.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 26500
$ws.Range("I14").Value = 50000
$ws.Range("J14").Value = 3000
$ws.Range("K14").Value = 50000
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = -49828
$ws.Range("N14").Value = -3344
# Row 22
$ws.Range("H22").Value = 8069.5713
$ws.Range("I22").Value = 511.42856
$ws.Range("J22").Value = 15627.714
$ws.Range("K22").Value = 511.42856
$ws.Range("L22").Value = 15627.714
$ws.Range("M22").Value = -216.42856
$ws.Range("N22").Value = -16217.714
# Row 27
$ws.Range("H27").Value = 8069.5713
$ws.Range("I27").Value = 511.42856
$ws.Range("J27").Value = 15627.714
$ws.Range("K27").Value = 511.42856
$ws.Range("L27").Value = 15627.714
$ws.Range("M27").Value = -404.42856
$ws.Range("N27").Value = -15841.714
# Row 132
$ws.Range("H132").Value = 4081.6667
$ws.Range("I132").Value = 2435.25
$ws.Range("J132").Value = 7374.5
$ws.Range("K132").Value = 7305.75
$ws.Range("L132").Value = 22123.5
$ws.Range("M132").Value = -4775.75
$ws.Range("N132").Value = -27183.5
# Row 136
$ws.Range("H136").Value = 4879.391
$ws.Range("I136").Value = 2036.5
$ws.Range("J136").Value = 11377.429
$ws.Range("K136").Value = 6109.5
$ws.Range("L136").Value = 34132.287
$ws.Range("M136").Value = -3559.5
$ws.Range("N136").Value = -39232.287

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 100006600
$ws.Range("I132").Value = 125007000
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 375021000
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -375018470
$ws.Range("N132").Value = -20057
# Row 136
$ws.Range("H136").Value = 11941907
$ws.Range("I136").Value = 16717681
$ws.Range("J136").Value = 2472.5
$ws.Range("K136").Value = 50153043
$ws.Range("L136").Value = 7417.5
$ws.Range("M136").Value = -50150493
$ws.Range("N136").Value = -12517.5
